# The sheet originally holds a single numeric row (A1:E1 = 100,105,110,115,120).
# The new version replaces that with a 2-row table of text-typed numbers:
#   row 1: 1  2  3  4  5
#   row 2: 60 90 30 120 60
# Values are entered as genuine text (shared strings), not numbers, so a
# leading apostrophe is used to force text entry (mirrors how Excel itself
# stores a quote-prefixed numeric-looking entry). Cells are written column
# by column (A1,A2,B1,B2,...) so the shared-string table is built up in
# that interleaved order. Formatting is cleared afterwards so the cells
# keep the default (unstyled) appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "'1"
$ws.Range("A2").Value = "'60"

$ws.Range("B1").Value = "'2"
$ws.Range("B2").Value = "'90"

$ws.Range("C1").Value = "'3"
$ws.Range("C2").Value = "'30"

$ws.Range("D1").Value = "'4"
$ws.Range("D2").Value = "'120"

$ws.Range("E1").Value = "'5"
$ws.Range("E2").Value = "'60"

$ws.Range("A1:E2").ClearFormats()
